# Generate Report for Handoff
# The da5e691d-f6b3-4ecb-a132-7b63c4c63144.md file moves from
# "Handed back: in sync with en-US" to "Ready for handoff" status, with
# updated handoff timestamps and (for the two locale sheets) a new
# "version not latest" error message in the Error Detail column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row for da5e691d-...md (row 3)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-06 08:39:29"

# ---------------------------------------------------------------------
# zh-cn sheet: row for da5e691d-...md (row 3)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-09-06 08:39:17"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/67cb649e67b0d62b4709669f6aa2183e3dd3f34b/e2e/da5e691d-f6b3-4ecb-a132-7b63c4c63144.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41768a9770881096595cb40dce5ac127ffe592ad/e2e/da5e691d-f6b3-4ecb-a132-7b63c4c63144.md."
$wsZhCn.Range("P:P").ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------
# de-de sheet: row for da5e691d-...md (row 3)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-09-06 08:39:29"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/67cb649e67b0d62b4709669f6aa2183e3dd3f34b/e2e/da5e691d-f6b3-4ecb-a132-7b63c4c63144.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41768a9770881096595cb40dce5ac127ffe592ad/e2e/da5e691d-f6b3-4ecb-a132-7b63c4c63144.md."
$wsDeDe.Range("P:P").ColumnWidth = 39.1666666666667
